$wb = $excel.ActiveWorkbook

$wsDomains = $wb.Worksheets.Item("Domains")
$wsCells   = $wb.Worksheets.Item("Cells")
$wsEdges   = $wb.Worksheets.Item("Edges")
$wsModules = $wb.Worksheets.Item("Modules")

# --- Rename header row labels to lowercase / snake_case "query param" style keys ---

# Cells sheet
$wsCells.Range("A1").Value = "name"
$wsCells.Range("B1").Value = "domain"

# Edges sheet
$wsEdges.Range("A1").Value = "from"
$wsEdges.Range("B1").Value = "to"

# Modules sheet
$wsModules.Range("A1").Value = "number"
$wsModules.Range("B1").Value = "title"
$wsModules.Range("C1").Value = "cell_1"
$wsModules.Range("D1").Value = "cell_2"
$wsModules.Range("E1").Value = "cell_3"
$wsModules.Range("F1").Value = "cell_4"
$wsModules.Range("G1").Value = "cell_5"
$wsModules.Range("H1").Value = "cell_6"
$wsModules.Range("I1").Value = "cell_7"
$wsModules.Range("J1").Value = "cell_8"
$wsModules.Range("K1").Value = "cell_9"
$wsModules.Range("L1").Value = "cell_10"
$wsModules.Range("M1").Value = "cell_11"

# --- New "query param" column (N) on the Modules sheet, plus column widths ---
$wsModules.Range("N1").Font.Bold = $true
$wsModules.Columns.Item(2).ColumnWidth = 11
$wsModules.Columns.Item(14).ColumnWidth = 8.5

# --- View / selection state ---

$wsDomains.Activate()
$wsDomains.Range("A3").Select()

$wsCells.Activate()
$wsCells.Application.ActiveWindow.ScrollRow = 1
$wsCells.Range("B5").Select()

$wsEdges.Activate()
$wsEdges.Range("B2").Select()

$wsModules.Activate()
$wsModules.Range("D9").Select()
